# Apply fixture update to the Man City schedule sheet:
#  - remove the "Manchester City v Wolverhampton Wanderers" fixture (row 1)
#  - insert a new "Arsenal v Manchester City" fixture before the Chelsea fixture

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the Wolverhampton Wanderers row entirely (row 1); this shifts every
# following row up by one.
$ws.Rows.Item(1).Delete()

# After the deletion, the old "26/12/2021" Leicester City row is now row 2,
# and the Chelsea row is now row 3. Insert a fresh row before Chelsea and
# populate it with the new Arsenal fixture.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "Arsenal v Manchester City "
$ws.Cells.Item(3, 2).Value = "01/01/2022 12:30 | Premier League"
